# Natmi following Dr Hou advice:
# The ligand/receptor (Tnfsf13/Fas) cluster-pair table is recomputed with a
# third sending/target cluster (FAPs) added alongside the existing ECs/sCs
# clusters, so the 2x3 sender/target grid becomes a full 3x3 grid (9 rows)
# and every statistic column is refreshed with the newly computed values.
$data = @(
  @("ECs","Tnfsf13","Fas","ECs",3,1,0.291812,0.875436,0.1179458223311005,0.1179458223311005,3,1,15.414045,46.242135,0.5885168832403832,0.5885168832403832,4.498003299540001,40.48202969586001,0.06941310774952325,0.06941310774952325),
  @("ECs","Tnfsf13","Fas","FAPs",3,1,0.291812,0.875436,0.1179458223311005,0.1179458223311005,3,1,6.673535333333334,20.020606,0.254799321953965,0.254799321953965,1.947417692690667,17.526759234216,0.03005251555726723,0.03005251555726723),
  @("ECs","Tnfsf13","Fas","sCs",3,1,0.291812,0.875436,0.1179458223311005,0.1179458223311005,3,1,4.103758333333333,12.311275,0.1566837948056518,0.1566837948056518,1.197525926766666,10.7777333409,0.01848019902431001,0.01848019902431001),
  @("FAPs","Tnfsf13","Fas","ECs",2,0.6666666666666666,0.514066,1.542198,0.2077773947009016,0.2077773947009016,3,1,15.414045,46.242135,0.5885168832403832,0.5885168832403832,7.923836456970001,71.31452811273,0.1222805047371815,0.1222805047371815),
  @("FAPs","Tnfsf13","Fas","FAPs",2,0.6666666666666666,0.514066,1.542198,0.2077773947009016,0.2077773947009016,3,1,6.673535333333334,20.020606,0.254799321953965,0.254799321953965,3.430637614665334,30.875738531988,0.0529415392871511,0.0529415392871511),
  @("FAPs","Tnfsf13","Fas","sCs",2,0.6666666666666666,0.514066,1.542198,0.2077773947009016,0.2077773947009016,3,1,4.103758333333333,12.311275,0.1566837948056518,0.1566837948056518,2.109602631383333,18.98642368245,0.03255535067656899,0.03255535067656899),
  @("sCs","Tnfsf13","Fas","ECs",3,1,1.668241,5.004723,0.674276782967998,0.674276782967998,3,1,15.414045,46.242135,0.5885168832403832,0.5885168832403832,25.71434184484501,231.429076603605,0.3968232707536785,0.3968232707536785),
  @("sCs","Tnfsf13","Fas","FAPs",3,1,1.668241,5.004723,0.674276782967998,0.674276782967998,3,1,6.673535333333334,20.020606,0.254799321953965,0.254799321953965,11.13306525801533,100.197587322138,0.1718052671095467,0.1718052671095467),
  @("sCs","Tnfsf13","Fas","sCs",3,1,1.668241,5.004723,0.674276782967998,0.674276782967998,3,1,4.103758333333333,12.311275,0.1566837948056518,0.1566837948056518,6.846057905758332,61.614521151825,0.1056482451047728,0.1056482451047728),
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
